$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 4113.1665
$ws.Range("J58").Value = 4736
$ws.Range("L58").Value = 14208
$ws.Range("N58").Value = -14508
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("M111").ClearContents()
$ws.Range("H137").Value = 3695.75
$ws.Range("I137").Value = 2392.5
$ws.Range("K137").Value = 7177.5
$ws.Range("M137").Value = -4627.5
$ws.Range("H138").Value = 11033.424
$ws.Range("J138").Value = 11121.876
$ws.Range("L138").Value = 33365.628
$ws.Range("N138").Value = -43645.628
$ws.Range("H141").Value = 3874.5
$ws.Range("I141").Value = 3874.5
$ws.Range("K141").Value = 11623.5
$ws.Range("M141").Value = -6443.5

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 75.8421
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H32").Value = 7776.6665
$ws.Range("I32").Value = 7776.6665
$ws.Range("K32").Value = 7776.6665
$ws.Range("M32").Value = -7489.6665
$ws.Range("H45").Value = 7116.16
$ws.Range("I45").Value = 14947.4
$ws.Range("K45").Value = 14947.4
$ws.Range("M45").Value = -14570.4
$ws.Range("H61").Value = 6379.3125
$ws.Range("I61").Value = 5760.1665
$ws.Range("J61").Value = 8236.75
$ws.Range("K61").Value = 5760.1665
$ws.Range("L61").Value = 8236.75
$ws.Range("M61").Value = -5548.1665
$ws.Range("N61").Value = -8660.75
$ws.Range("H74").Value = 5165.6665
$ws.Range("I74").Value = 4950.1
$ws.Range("K74").Value = 4950.1
$ws.Range("M74").Value = -4076.1
$ws.Range("H77").Value = 5165.6665
$ws.Range("I77").Value = 4950.1
$ws.Range("K77").Value = 24750.5
$ws.Range("M77").Value = -20382.5
$ws.Range("H122").Value = 13610.667
$ws.Range("I122").Value = 15566.846
$ws.Range("K122").Value = 46700.538
$ws.Range("M122").Value = -44250.538
$ws.Range("H132").Value = 5942.6924
$ws.Range("I132").Value = 6404.3
$ws.Range("K132").Value = 19212.9
$ws.Range("M132").Value = -16682.9
$ws.Range("H136").Value = 6379.3125
$ws.Range("I136").Value = 5760.1665
$ws.Range("J136").Value = 8236.75
$ws.Range("K136").Value = 17280.4995
$ws.Range("L136").Value = 24710.25
$ws.Range("M136").Value = -14730.4995
$ws.Range("N136").Value = -29810.25

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 75.8421
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H34").Value = 14300
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 14300
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 14300
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -14528
$ws.Range("H94").Value = 1450.3334
$ws.Range("J94").Value = 826.6667
$ws.Range("L94").Value = 826.6667
$ws.Range("N94").Value = -1728.6667
$ws.Range("H107").Value = 1245.2609
$ws.Range("I107").Value = 1199.2106
$ws.Range("J107").Value = 1464
$ws.Range("K107").Value = 1199.2106
$ws.Range("L107").Value = 1464
$ws.Range("M107").Value = 720.7893999999999
$ws.Range("N107").Value = -5304

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6659.143
$ws.Range("I31").Value = 5819.5
$ws.Range("K31").Value = 5819.5
$ws.Range("M31").Value = -5524.5
$ws.Range("H34").Value = 6659.143
$ws.Range("I34").Value = 5819.5
$ws.Range("K34").Value = 5819.5
$ws.Range("M34").Value = -5617.5
$ws.Range("H58").Value = 4764.56
$ws.Range("I58").Value = 4701.1665
$ws.Range("K58").Value = 4701.1665
$ws.Range("M58").Value = -4498.1665
$ws.Range("H62").Value = 1874.12
$ws.Range("J62").Value = 2333.3333
$ws.Range("L62").Value = 2333.3333
$ws.Range("N62").Value = -3581.3333
$ws.Range("H65").Value = 1874.12
$ws.Range("J65").Value = 2333.3333
$ws.Range("L65").Value = 11666.6665
$ws.Range("N65").Value = -17906.6665
$ws.Range("H132").Value = 3308.353
$ws.Range("I132").Value = 2913.0833
$ws.Range("J132").Value = 4257
$ws.Range("K132").Value = 8739.249899999999
$ws.Range("L132").Value = 12771
$ws.Range("M132").Value = -6209.249899999999
$ws.Range("N132").Value = -17831
$ws.Range("H136").Value = 4764.56
$ws.Range("I136").Value = 4701.1665
$ws.Range("K136").Value = 14103.4995
$ws.Range("M136").Value = -11553.4995

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 5360.375
$ws.Range("J94").Value = 6387.6
$ws.Range("L94").Value = 19162.8
$ws.Range("N94").Value = -20514.8
$ws.Range("H122").Value = 2497.5
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 2497.5
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 22477.5
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -27377.5
$ws.Range("H129").Value = 2549
$ws.Range("J129").Value = 2549
$ws.Range("L129").Value = 7647
$ws.Range("N129").Value = -17647
$ws.Range("H131").Value = 4019.5715
$ws.Range("I131").Value = 2332.3333
$ws.Range("J131").Value = 5285
$ws.Range("K131").Value = 6996.999899999999
$ws.Range("L131").Value = 15855
$ws.Range("M131").Value = -1956.999899999999
$ws.Range("N131").Value = -25935
$ws.Range("H137").Value = 8346.286
$ws.Range("I137").Value = 6714.5
$ws.Range("J137").Value = 8999
$ws.Range("K137").Value = 20143.5
$ws.Range("L137").Value = 26997
$ws.Range("M137").Value = -15043.5
$ws.Range("N137").Value = -37197

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 89352.92999999999
$ws.Range("J141").Value = 89352.92999999999
$ws.Range("L141").Value = 89352.92999999999
$ws.Range("N141").Value = -99712.92999999999

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 15184.385
$ws.Range("I40").Value = 14665.223
$ws.Range("K40").Value = 14665.223
$ws.Range("M40").Value = -14529.223
$ws.Range("H46").Value = 2117.5
$ws.Range("I46").Value = 1500.25
$ws.Range("K46").Value = 1500.25
$ws.Range("M46").Value = -1312.25
$ws.Range("H68").Value = 3833.5
$ws.Range("I68").Value = 3833.5
$ws.Range("K68").Value = 3833.5
$ws.Range("M68").Value = -3084.5
$ws.Range("H71").Value = 3833.5
$ws.Range("I71").Value = 3833.5
$ws.Range("K71").Value = 19167.5
$ws.Range("M71").Value = -15423.5
$ws.Range("H100").Value = 4445
$ws.Range("I100").Value = 4281.1665
$ws.Range("K100").Value = 4281.1665
$ws.Range("M100").Value = -3740.1665
$ws.Range("H122").Value = 5315.4
$ws.Range("I122").Value = 5139.5757
$ws.Range("J122").Value = 6144.2856
$ws.Range("K122").Value = 15418.7271
$ws.Range("L122").Value = 18432.8568
$ws.Range("M122").Value = -12968.7271
$ws.Range("N122").Value = -23332.8568
$ws.Range("H132").Value = 23679.1
$ws.Range("I132").Value = 15446
$ws.Range("K132").Value = 46338
$ws.Range("M132").Value = -43808
$ws.Range("H136").Value = 8327.741
$ws.Range("I136").Value = 6996.2856
$ws.Range("K136").Value = 20988.8568
$ws.Range("M136").Value = -18438.8568

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3012.8
$ws.Range("I81").Value = 3012.8
$ws.Range("K81").Value = 6025.6
$ws.Range("M81").Value = -4964.6
$ws.Range("H84").Value = 3012.8
$ws.Range("I84").Value = 3012.8
$ws.Range("K84").Value = 30128
$ws.Range("M84").Value = -24824
$ws.Range("H122").Value = 2272.6667
$ws.Range("I122").Value = 2272.6667
$ws.Range("K122").Value = 6818.000100000001
$ws.Range("M122").Value = -4368.000100000001
$ws.Range("H126").Value = 1839.7407
$ws.Range("I126").Value = 1679.7307
$ws.Range("K126").Value = 5039.1921
$ws.Range("M126").Value = -2569.1921
$ws.Range("H132").Value = 5988.88
$ws.Range("I132").Value = 5410
$ws.Range("K132").Value = 16230
$ws.Range("M132").Value = -13700
